$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.059.63'
$ws.Range("E2").Value = '  +2.46%  '

# Row 3
$ws.Range("D3").Value = '3.603.50'
$ws.Range("E3").Value = '  +1.07%  '

# Row 4
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '204.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.67%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '564.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.80%  '

# Row 7
$ws.Range("D7").Value = '3.596.38'
$ws.Range("E7").Value = '  +1.04%  '

# Row 8
$ws.Range("E8").Value = '  +1.26%  '

# Row 9
$ws.Range("E9").Value = '  -0.17%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '61.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +15.65%  '

# Row 12
$ws.Range("E12").Value = '  +3.44%  '

# Row 13
$ws.Range("E13").Value = '  +10.00%  '

# Row 14
$ws.Range("E14").Value = '  +2.11%  '

# Row 15
$ws.Range("D15").Value = '4.191.68'
$ws.Range("E15").Value = '  +1.39%  '

# Row 16
$ws.Range("D16").Value = '3.608.55'
$ws.Range("E16").Value = '  +1.39%  '

# Row 17
$ws.Range("E17").Value = '  +0.66%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.15%  '

# Row 19
$ws.Range("D19").Value = '67.918.37'
$ws.Range("E19").Value = '  +2.42%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.02%  '

# Row 21
$ws.Range("E21").Value = '  +1.90%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '402.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.56%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.07%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.51%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.74%  '

# Row 26
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.06%  '

# Row 27
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.96'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.98%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.62'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.13'
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +17.17%  '

# Row 31
$ws.Range("E31").Value = '  +4.93%  '

# Row 32
$ws.Range("E32").Value = '  +1.59%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '673.05'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.55%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.22'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.33%  '

# Row 35
$ws.Range("E35").Value = '  +0.64%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '63.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.04%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.36%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.421'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.89%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.20%  '

# Row 40
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '3.307.99'
$ws.Range("E40").Value = '  +9.93%  '

# Row 41
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0769'
$ws.Range("E41").Value = '  +0.26%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.70%  '

# Row 43
$ws.Range("E43").Value = '  +3.85%  '

# Row 44
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.07'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +32.16%  '

# Row 45
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.78'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.68%  '

# Row 46
$ws.Range("E46").Value = '  -0.15%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0419'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.30%  '

# Row 48
$ws.Range("E48").Value = '  +11.83%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.38%  '

# Row 50
$ws.Range("E50").Value = '  +0.44%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.44%  '
